$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Item 1)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "HEPARINA SODICA SUBCUT 5000UI"
$ws.Cells.Item(2, 3).Value = "5000ui"
$ws.Cells.Item(2, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(2, 5).Value = "'102980371"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 6).Value = "OK"

# Row 3 (Item 4)
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "PROMETAZINA 25MG, CLORIDRATO"
$ws.Cells.Item(3, 3).Value = "25mg"
$ws.Cells.Item(3, 4).Value = "LABORATÓRIO TEUTO BRASILEIRO S/A"
$ws.Cells.Item(3, 5).Value = "'103700321"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 6).Value = "OK"

# Row 4 (Item 5)
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "HALOPERIDOL 5MG"
$ws.Cells.Item(4, 3).Value = "5mg"
$ws.Cells.Item(4, 4).Value = "CELLERA FARMACÊUTICA S.A."
$ws.Cells.Item(4, 5).Value = "'112360011"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 6).Value = "OK"

# Row 5 (Item 6)
$ws.Cells.Item(5, 1).Value = 6
$ws.Cells.Item(5, 2).Value = "CLORPROMAZINA 40MG/ML SOL ORAL"
$ws.Cells.Item(5, 3).Value = "40mg/ml"
$ws.Cells.Item(5, 4).Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Cells.Item(5, 5).Value = "Último registro encontrado: 183260385"
$ws.Cells.Item(5, 6).Value = "Pendente"

# Row 6 (Item 7)
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = "HALOPERIDOL 2MG/ML SOL ORAL"
$ws.Cells.Item(6, 3).Value = "2mg/ml"
$ws.Cells.Item(6, 4).Value = "CELLERA FARMACÊUTICA S.A."
$ws.Cells.Item(6, 5).Value = "'112360011"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 6).Value = "OK"

# Row 7 (Item 8)
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "CLORPROMAZINA 25MG"
$ws.Cells.Item(7, 3).Value = "25mg"
$ws.Cells.Item(7, 4).Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Cells.Item(7, 5).Value = "'183260385"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 6).Value = "OK"

# Row 8 (Item 9)
$ws.Cells.Item(8, 1).Value = 9
$ws.Cells.Item(8, 2).Value = "CODEINA 30MG"
$ws.Cells.Item(8, 3).Value = "30mg"
$ws.Cells.Item(8, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(8, 5).Value = "'102980199"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 6).Value = "OK"

# Row 9 (Item 10)
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "IMIPRAMINA 25MG"
$ws.Cells.Item(9, 3).Value = "25mg"
$ws.Cells.Item(9, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(9, 5).Value = "'102980023"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 6).Value = "OK"

# Row 10 (Item 11)
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "RISPERIDONA 3MG"
$ws.Cells.Item(10, 3).Value = "3mg"
$ws.Cells.Item(10, 4).Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Cells.Item(10, 5).Value = "'103920197"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 6).Value = "OK"

# Row 11 (Item 13)
$ws.Cells.Item(11, 1).Value = 13
$ws.Cells.Item(11, 2).Value = "RISPERIDONA 1MG"
$ws.Cells.Item(11, 3).Value = "1mg"
$ws.Cells.Item(11, 4).Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Cells.Item(11, 5).Value = "'103920197"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 6).Value = "OK"

# Row 12 (Item 14)
$ws.Cells.Item(12, 1).Value = 14
$ws.Cells.Item(12, 2).Value = "LEVOMEPROMAZINA 4% GOTAS"
$ws.Cells.Item(12, 3).Value = "'4%"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Cells.Item(12, 5).Value = "Último registro encontrado: 183260316"
$ws.Cells.Item(12, 6).Value = "Pendente"

# Row 13 (Item 16)
$ws.Cells.Item(13, 1).Value = 16
$ws.Cells.Item(13, 2).Value = "LIDOCAINA 2% C/ VASO CONSTRITO"
$ws.Cells.Item(13, 3).Value = "'2%"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(13, 5).Value = "'102980249"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 6).Value = "OK"

# Row 14 (Item 17)
$ws.Cells.Item(14, 1).Value = 17
$ws.Cells.Item(14, 2).Value = "NITRATO DE CERIO +SULFADIAZINA"
$ws.Cells.Item(14, 3).Value = "Concentração não encontrada"
$ws.Cells.Item(14, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(14, 5).Value = "Último registro encontrado: 102980560"
$ws.Cells.Item(14, 6).Value = "Pendente"

# Row 15 (Item 18)
$ws.Cells.Item(15, 1).Value = 18
$ws.Cells.Item(15, 2).Value = "COLAGENASE+CLORAFENICOL POMADA 30g"
$ws.Cells.Item(15, 3).Value = "30g"
$ws.Cells.Item(15, 4).Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Cells.Item(15, 5).Value = "'102980431"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 6).Value = "OK"
